# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet, with the
#    per-fund holdings table for the new quarter.
# 2) Prepend a "2022-Q1" row to the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")   # same A1:H* layout/styling to copy from
$totalSheet = $wb.Worksheets.Item("总计")

# Snapshot the existing "总计" data rows (2021-Q4 .. 2020-Q4) before touching
# anything, so they can be rewritten after the sheet is recreated below.
$existingRows = @()
for ($r = 2; $r -le 6; $r++) {
    $existingRows += ,@(
        $totalSheet.Cells.Item($r, 2).Value2,
        $totalSheet.Cells.Item($r, 3).Value2,
        $totalSheet.Cells.Item($r, 4).Value2
    )
}

# Recreate the "总计" sheet (delete + re-add at the end) so that the brand
# new "2022-Q1" sheet inherits the lower sheetId (6) and "总计" gets the
# fresh one (7) - matching the order the sheets were authored in.
$totalSheet.Delete()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Add($null, $newSheet)
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet
# ---------------------------------------------------------------------------

# Copy header row (row 1, columns B..H only - A1 has no cell in the source)
# + first data row (row 2) so we inherit the exact fonts/borders (bold +
# bordered column A / header, plain data cells).
$template.Range("B1:H2").Copy($newSheet.Range("B1"))
# Stretch the row-2 style pattern down across all 11 data rows (2..12).
$template.Range("A2:H2").Copy($newSheet.Range("A2:A12"))

# Columns B..G hold numeric-looking values that must be stored as TEXT
# (leading zeros in fund codes, trailing zeros in percentages, etc.), matching
# the source data. Force text format before writing, then strip the
# formatting change back out so no stray style index is left on the cells.
$textRange = $newSheet.Range("B2:G12")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "001120"
$newSheet.Range("C2").Value = "东方睿鑫热点挖掘灵活配置混合A"
$newSheet.Range("D2").Value = "1.74"
$newSheet.Range("E2").Value = "78.43"
$newSheet.Range("F2").Value = "4.90"
$newSheet.Range("G2").Value = "0.0853"

$newSheet.Range("B3").Value = "008895"
$newSheet.Range("C3").Value = "申万菱信量化对冲策略灵活配置混合"
$newSheet.Range("D3").Value = "8.57"
$newSheet.Range("E3").Value = "42.67"
$newSheet.Range("F3").Value = "0.95"
$newSheet.Range("G3").Value = "0.0814"

$newSheet.Range("B4").Value = "008997"
$newSheet.Range("C4").Value = "同泰竞争优势混合A"
$newSheet.Range("D4").Value = "1.33"
$newSheet.Range("E4").Value = "94.35"
$newSheet.Range("F4").Value = "5.53"
$newSheet.Range("G4").Value = "0.0735"

$newSheet.Range("B5").Value = "012496"
$newSheet.Range("C5").Value = "同泰行业优选股票A"
$newSheet.Range("D5").Value = "1.10"
$newSheet.Range("E5").Value = "94.13"
$newSheet.Range("F5").Value = "5.24"
$newSheet.Range("G5").Value = "0.0576"

$newSheet.Range("B6").Value = "004205"
$newSheet.Range("C6").Value = "东方支柱产业灵活配置混合"
$newSheet.Range("D6").Value = "0.96"
$newSheet.Range("E6").Value = "84.14"
$newSheet.Range("F6").Value = "5.64"
$newSheet.Range("G6").Value = "0.0541"

$newSheet.Range("B7").Value = "001121"
$newSheet.Range("C7").Value = "东方睿鑫热点挖掘灵活配置混合C"
$newSheet.Range("D7").Value = "1.10"
$newSheet.Range("E7").Value = "78.43"
$newSheet.Range("F7").Value = "4.90"
$newSheet.Range("G7").Value = "0.0539"

$newSheet.Range("B8").Value = "008998"
$newSheet.Range("C8").Value = "同泰竞争优势混合C"
$newSheet.Range("D8").Value = "0.95"
$newSheet.Range("E8").Value = "94.35"
$newSheet.Range("F8").Value = "5.53"
$newSheet.Range("G8").Value = "0.0525"

$newSheet.Range("B9").Value = "004244"
$newSheet.Range("C9").Value = "东方周期优选灵活配置混合"
$newSheet.Range("D9").Value = "0.62"
$newSheet.Range("E9").Value = "84.13"
$newSheet.Range("F9").Value = "5.79"
$newSheet.Range("G9").Value = "0.0359"

$newSheet.Range("B10").Value = "012497"
$newSheet.Range("C10").Value = "同泰行业优选股票C"
$newSheet.Range("D10").Value = "0.28"
$newSheet.Range("E10").Value = "94.13"
$newSheet.Range("F10").Value = "5.24"
$newSheet.Range("G10").Value = "0.0147"

$newSheet.Range("B11").Value = "005443"
$newSheet.Range("C11").Value = "国金量化多策略灵活配置混合"
$newSheet.Range("D11").Value = "0.51"
$newSheet.Range("E11").Value = "64.10"
$newSheet.Range("F11").Value = "0.92"
$newSheet.Range("G11").Value = "0.0047"

$newSheet.Range("B12").Value = "006195"
$newSheet.Range("C12").Value = "国金量化多因子股票"
$newSheet.Range("D12").Value = "0.09"
$newSheet.Range("E12").Value = "80.71"
$newSheet.Range("F12").Value = "0.89"
$newSheet.Range("G12").Value = "0.0008"

$textRange.ClearFormats()

# Column H (仓位排名) is a genuine number, and column A is the 0-based index.
$newSheet.Range("A2").Value = 0
$newSheet.Range("H2").Value = 4
$newSheet.Range("A3").Value = 1
$newSheet.Range("H3").Value = 10
$newSheet.Range("A4").Value = 2
$newSheet.Range("H4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("H5").Value = 5
$newSheet.Range("A6").Value = 4
$newSheet.Range("H6").Value = 6
$newSheet.Range("A7").Value = 5
$newSheet.Range("H7").Value = 4
$newSheet.Range("A8").Value = 6
$newSheet.Range("H8").Value = 2
$newSheet.Range("A9").Value = 7
$newSheet.Range("H9").Value = 2
$newSheet.Range("A10").Value = 8
$newSheet.Range("H10").Value = 5
$newSheet.Range("A11").Value = 9
$newSheet.Range("H11").Value = 3
$newSheet.Range("A12").Value = 10
$newSheet.Range("H12").Value = 6

# Re-assert the header text (already correct via the copy, but explicit is safe)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 2. Rebuild the "总计" sheet: header, new 2022-Q1 row, then the rows that
#    already existed, each shifted down by one.
# ---------------------------------------------------------------------------

# Borrow the same bold+bordered header / column-A style from the template
# sheet (columns A..D line up with the same style pattern used here; A1 has
# no cell in the source so copy B1:D1 only to avoid materialising a stray
# empty A1 cell).
$template.Range("B1:D2").Copy($totalSheet.Range("B1"))
$template.Range("A2:D2").Copy($totalSheet.Range("A2:A7"))

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 0.51

for ($i = 0; $i -lt $existingRows.Count; $i++) {
    $r = $i + 3
    $totalSheet.Cells.Item($r, 1).Value = $i + 1
    $totalSheet.Cells.Item($r, 2).Value = $existingRows[$i][0]
    $totalSheet.Cells.Item($r, 3).Value = $existingRows[$i][1]
    $totalSheet.Cells.Item($r, 4).Value = $existingRows[$i][2]
}
